# Scheduled price-refresh update across all leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 451.72726
$ws.Range("I41").Value = 43.666668
$ws.Range("K41").Value = 43.666668
$ws.Range("M41").Value = 396.333332

$ws.Range("H55").Value = 233.28572
$ws.Range("I55").Value = 77.125
$ws.Range("J55").Value = 441.5
$ws.Range("K55").Value = 77.125
$ws.Range("L55").Value = 441.5
$ws.Range("M55").Value = 136.875
$ws.Range("N55").Value = -869.5

$ws.Range("H113").Value = 3499.5
$ws.Range("I113").Value = 3499.5
$ws.Range("K113").Value = 3499.5
$ws.Range("M113").Value = -245.5

$ws.Range("H137").Value = 3046.8235
$ws.Range("I137").Value = 1391.375
$ws.Range("J137").Value = 4518.3335
$ws.Range("K137").Value = 4174.125
$ws.Range("L137").Value = 13555.0005
$ws.Range("M137").Value = -1624.125
$ws.Range("N137").Value = -18655.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 196.2
$ws.Range("I5").Value = 195.5
$ws.Range("J5").Value = 199
$ws.Range("K5").Value = 195.5
$ws.Range("L5").Value = 199
$ws.Range("M5").Value = -83.5
$ws.Range("N5").Value = -423

$ws.Range("H61").Value = 13499.5
$ws.Range("I61").Value = 19999
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 19999
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -19787
$ws.Range("N61").Value = -7424

$ws.Range("H74").Value = 1797.9166
$ws.Range("I74").Value = 1790.1111
$ws.Range("J74").Value = 1821.3334
$ws.Range("K74").Value = 1790.1111
$ws.Range("L74").Value = 1821.3334
$ws.Range("M74").Value = -916.1111000000001
$ws.Range("N74").Value = -3569.3334

$ws.Range("H77").Value = 1797.9166
$ws.Range("I77").Value = 1790.1111
$ws.Range("J77").Value = 1821.3334
$ws.Range("K77").Value = 8950.5555
$ws.Range("L77").Value = 9106.666999999999
$ws.Range("M77").Value = -4582.5555
$ws.Range("N77").Value = -17842.667

$ws.Range("H122").Value = 2256.5
$ws.Range("I122").Value = 2256.5
$ws.Range("K122").Value = 6769.5
$ws.Range("M122").Value = -4319.5

$ws.Range("H132").Value = 5998.5
$ws.Range("I132").Value = 5998
$ws.Range("K132").Value = 17994
$ws.Range("M132").Value = -15464

$ws.Range("H136").Value = 13499.5
$ws.Range("I136").Value = 19999
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 59997
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -57447
$ws.Range("N136").Value = -26100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 196.2
$ws.Range("I4").Value = 195.5
$ws.Range("J4").Value = 199
$ws.Range("K4").Value = 195.5
$ws.Range("L4").Value = 199
$ws.Range("M4").Value = -80.5
$ws.Range("N4").Value = -429

$ws.Range("H22").Value = 648.2857
$ws.Range("I22").Value = 656.3333
$ws.Range("K22").Value = 656.3333
$ws.Range("M22").Value = -483.3333

$ws.Range("H86").Value = 1762
$ws.Range("I86").Value = 1615.1428
$ws.Range("K86").Value = 1615.1428
$ws.Range("M86").Value = -492.1428000000001

$ws.Range("H89").Value = 1762
$ws.Range("I89").Value = 1615.1428
$ws.Range("K89").Value = 8075.714
$ws.Range("M89").Value = -2459.714

$ws.Range("H94").Value = 1500
$ws.Range("I94").Value = 1500
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1500
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1049
$ws.Range("N94").ClearContents()

$ws.Range("H107").Value = 2807.889
$ws.Range("I107").Value = 2807.889
$ws.Range("K107").Value = 2807.889
$ws.Range("M107").Value = -887.8890000000001

$ws.Range("H129").Value = 80000
$ws.Range("J129").Value = 80000
$ws.Range("L129").Value = 80000
$ws.Range("N129").Value = -90000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1046.2222
$ws.Range("I22").Value = 673.6667
$ws.Range("K22").Value = 673.6667
$ws.Range("M22").Value = -323.6667

$ws.Range("H58").Value = 9685.637000000001
$ws.Range("I58").Value = 9566.625
$ws.Range("K58").Value = 9566.625
$ws.Range("M58").Value = -9363.625

$ws.Range("H99").Value = 6597.737
$ws.Range("I99").Value = 5594.7144
$ws.Range("J99").Value = 9406.200000000001
$ws.Range("K99").Value = 5594.7144
$ws.Range("L99").Value = 9406.200000000001
$ws.Range("M99").Value = -4096.7144
$ws.Range("N99").Value = -12402.2

$ws.Range("H107").Value = 355.33334
$ws.Range("I107").Value = 205.64706
$ws.Range("K107").Value = 205.64706
$ws.Range("M107").Value = 1714.35294

$ws.Range("H126").Value = 6597.737
$ws.Range("I126").Value = 5594.7144
$ws.Range("J126").Value = 9406.200000000001
$ws.Range("K126").Value = 16784.1432
$ws.Range("L126").Value = 28218.6
$ws.Range("M126").Value = -14314.1432
$ws.Range("N126").Value = -33158.60000000001

$ws.Range("H132").Value = 11677.5
$ws.Range("J132").Value = 13199.25
$ws.Range("L132").Value = 39597.75
$ws.Range("N132").Value = -44657.75

$ws.Range("H136").Value = 9685.637000000001
$ws.Range("I136").Value = 9566.625
$ws.Range("K136").Value = 28699.875
$ws.Range("M136").Value = -26149.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 157221
$ws.Range("I2").Value = 32
$ws.Range("J2").Value = 220096.6
$ws.Range("K2").Value = 192
$ws.Range("L2").Value = 1320579.6
$ws.Range("M2").Value = -79
$ws.Range("N2").Value = -1320805.6

$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H87").Value = 616.3333
$ws.Range("I87").Value = 616.3333
$ws.Range("K87").Value = 1848.9999
$ws.Range("M87").Value = -600.9999

$ws.Range("H90").Value = 616.3333
$ws.Range("I90").Value = 616.3333
$ws.Range("K90").Value = 5546.9997
$ws.Range("M90").Value = 693.0002999999997

$ws.Range("H113").Value = 549.6667
$ws.Range("J113").Value = 499.5
$ws.Range("L113").Value = 1498.5
$ws.Range("N113").Value = -5838.5

$ws.Range("H114").Value = 1374.75
$ws.Range("I114").Value = 1250
$ws.Range("J114").Value = 1499.5
$ws.Range("K114").Value = 3750
$ws.Range("L114").Value = 4498.5
$ws.Range("M114").Value = -496
$ws.Range("N114").Value = -11006.5

$ws.Range("J117").Value = 450
$ws.Range("L117").Value = 1350
$ws.Range("N117").Value = -8234

$ws.Range("H131").Value = 1499.8572
$ws.Range("I131").Value = 900
$ws.Range("K131").Value = 2700
$ws.Range("M131").Value = 2340

$ws.Range("H137").Value = 1998
$ws.Range("J137").Value = 2031.1111
$ws.Range("L137").Value = 6093.3333
$ws.Range("N137").Value = -16293.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1263.6
$ws.Range("I97").Value = 1263.6
$ws.Range("K97").Value = 1263.6
$ws.Range("M97").Value = -767.5999999999999

$ws.Range("H132").Value = 2974.8333
$ws.Range("I132").Value = 2780
$ws.Range("K132").Value = 8340
$ws.Range("M132").Value = -5810

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 782.5454999999999
$ws.Range("I16").Value = 782.5454999999999
$ws.Range("K16").Value = 782.5454999999999
$ws.Range("M16").Value = -612.5454999999999

$ws.Range("H25").Value = 40000
$ws.Range("I25").Value = 40000
$ws.Range("K25").Value = 40000
$ws.Range("M25").Value = -39770

$ws.Range("H40").Value = 10833.333
$ws.Range("I40").Value = 4250
$ws.Range("K40").Value = 4250
$ws.Range("M40").Value = -4114

$ws.Range("H132").Value = 3919.7144
$ws.Range("J132").Value = 4566.6665
$ws.Range("L132").Value = 13699.9995
$ws.Range("N132").Value = -18759.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2233.3333
$ws.Range("I96").Value = 1680
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 1680
$ws.Range("L96").Value = 5000
$ws.Range("M96").Value = -307
$ws.Range("N96").Value = -7746

$ws.Range("H126").Value = 1244
$ws.Range("I126").Value = 1231.6666
$ws.Range("K126").Value = 3694.9998
$ws.Range("M126").Value = -1224.9998

$ws.Range("H136").Value = 3596.7334
$ws.Range("I136").Value = 3313.9092
$ws.Range("J136").Value = 4374.5
$ws.Range("K136").Value = 9941.7276
$ws.Range("L136").Value = 13123.5
$ws.Range("M136").Value = -7391.7276
$ws.Range("N136").Value = -18223.5
